# Updated cryptos list on Sun Aug 20 10:39:12 UTC 2023 with GitHub Actions
#
# Refreshes the "Price" (column D) and "Volume(1h)" (column E) figures for
# each coin row (rows 2-51) on the active worksheet with the latest scraped
# values. Column D values are stored as plain text (they come pre-formatted,
# e.g. "26.339.93" or "0.000008458") so numeric-looking values are written
# with a leading apostrophe to stop Excel from reinterpreting them as
# numbers/dates; the quote-prefix formatting that this introduces is then
# cleared by resetting the cell style back to "Normal" so the cell keeps its
# original (unstyled) appearance.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; Price = "26.339.93"; Volume = "  +1.14%  " },
    @{ Row = 3; Price = "1.683.32"; Volume = "  +0.95%  " },
    @{ Row = 4; Price = $null; Volume = "  +0.34%  " },
    @{ Row = 5; Price = "218.29"; Volume = "  +0.66%  " },
    @{ Row = 6; Price = "0.5523"; Volume = "  +8.31%  " },
    @{ Row = 7; Price = $null; Volume = "  +0.34%  " },
    @{ Row = 8; Price = "0.2702"; Volume = "  +1.75%  " },
    @{ Row = 9; Price = "0.06496"; Volume = "  +1.31%  " },
    @{ Row = 10; Price = "22.10"; Volume = "  +1.10%  " },
    @{ Row = 11; Price = "0.07558"; Volume = "  +1.51%  " },
    @{ Row = 12; Price = "4.547"; Volume = "  +0.88%  " },
    @{ Row = 13; Price = "1.678.74"; Volume = "  +0.67%  " },
    @{ Row = 14; Price = "0.5813"; Volume = "  -0.46%  " },
    @{ Row = 15; Price = "0.000008458"; Volume = "  -1.44%  " },
    @{ Row = 16; Price = "65.09"; Volume = "  +1.17%  " },
    @{ Row = 17; Price = "26.353.22"; Volume = "  +0.90%  " },
    @{ Row = 18; Price = "4.937"; Volume = "  -0.07%  " },
    @{ Row = 19; Price = $null; Volume = "  +0.34%  " },
    @{ Row = 20; Price = "10.92"; Volume = "  +1.28%  " },
    @{ Row = 21; Price = "191.22"; Volume = "  -0.34%  " },
    @{ Row = 22; Price = "6.236"; Volume = "  +0.43%  " },
    @{ Row = 23; Price = $null; Volume = "  +0.25%  " },
    @{ Row = 24; Price = "147.44"; Volume = "  +1.84%  " },
    @{ Row = 25; Price = "0.1321"; Volume = "  +10.37%  " },
    @{ Row = 26; Price = "7.899"; Volume = "  +3.59%  " },
    @{ Row = 27; Price = $null; Volume = "  +0.66%  " },
    @{ Row = 28; Price = "0.06339"; Volume = "  -2.30%  " },
    @{ Row = 29; Price = "1.392"; Volume = "  +4.88%  " },
    @{ Row = 30; Price = $null; Volume = "  +0.41%  " },
    @{ Row = 31; Price = "3.595"; Volume = "  +1.48%  " },
    @{ Row = 32; Price = "3.584"; Volume = "  +1.81%  " },
    @{ Row = 33; Price = "1.667"; Volume = "  +1.06%  " },
    @{ Row = 34; Price = "1.039"; Volume = "  +1.86%  " },
    @{ Row = 35; Price = "0.6210"; Volume = "  +1.58%  " },
    @{ Row = 36; Price = "2.401"; Volume = "  +1.39%  " },
    @{ Row = 37; Price = "2.721"; Volume = "  +1.44%  " },
    @{ Row = 38; Price = "6.236"; Volume = "  -0.54%  " },
    @{ Row = 39; Price = "1.112.32"; Volume = "  +1.93%  " },
    @{ Row = 40; Price = $null; Volume = "  +1.40%  " },
    @{ Row = 41; Price = "0.8714"; Volume = "  +0.91%  " },
    @{ Row = 42; Price = $null; Volume = "  +0.69%  " },
    @{ Row = 43; Price = "100.74"; Volume = $null },
    @{ Row = 44; Price = "1.831.78"; Volume = "  +0.87%  " },
    @{ Row = 45; Price = $null; Volume = "  -5.22%  " },
    @{ Row = 46; Price = $null; Volume = "  +1.58%  " },
    @{ Row = 47; Price = "8.181"; Volume = "  +1.64%  " },
    @{ Row = 48; Price = $null; Volume = "  -0.14%  " },
    @{ Row = 49; Price = "0.05273"; Volume = "  +0.81%  " },
    @{ Row = 50; Price = "0.4294"; Volume = "  +0.23%  " },
    @{ Row = 51; Price = "6.074"; Volume = "  +0.24%  " }
)

foreach ($update in $updates) {
    $row = $update.Row

    if ($null -ne $update.Price) {
        $priceCell = $ws.Range("D$row")
        $price = $update.Price

        # Coinranking renders big prices with "." as a thousands separator
        # (e.g. "26.339.93"), so only single-dot decimal values are at risk
        # of being reinterpreted by Excel as a real number.
        $looksNumeric = $price -match '^[0-9]+(\.[0-9]+)?$'

        if ($looksNumeric) {
            # Force text storage for values Excel would otherwise coerce to a
            # floating point number (e.g. "218.29", "0.000008458").
            $priceCell.Value = "'" + $price
            $priceCell.Style = "Normal"
        } else {
            $priceCell.Value = $price
        }
    }

    if ($null -ne $update.Volume) {
        $ws.Range("E$row").Value = $update.Volume
    }
}
